# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates to columns H-N across several sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 97: Materia Worth | Potent Spiritbond Potion
$ws.Range("H97").Value = 1762
$ws.Range("J97").Value = 1762
$ws.Range("L97").Value = 5286
$ws.Range("N97").Value = -6278

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 1615.3077
$ws.Range("I112").Value = 583.3333
$ws.Range("J112").Value = 1749.9131
$ws.Range("K112").Value = 1749.9999
$ws.Range("L112").Value = 5249.7393
$ws.Range("M112").Value = -641.9999
$ws.Range("N112").Value = -7465.7393

# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 505.76086
$ws.Range("I135").Value = 359.65115
$ws.Range("J135").Value = 2600
$ws.Range("K135").Value = 3236.86035
$ws.Range("L135").Value = 23400
$ws.Range("M135").Value = -701.8603499999999
$ws.Range("N135").Value = -28470

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 3497.818
$ws.Range("I137").Value = 3640.2173
$ws.Range("J137").Value = 3170.3
$ws.Range("K137").Value = 10920.6519
$ws.Range("L137").Value = 9510.900000000001
$ws.Range("M137").Value = -8370.651899999999
$ws.Range("N137").Value = -14610.9

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3208.36
$ws.Range("I138").Value = 1669.2559
$ws.Range("J138").Value = 5276.5312
$ws.Range("K138").Value = 5007.7677
$ws.Range("L138").Value = 15829.5936
$ws.Range("M138").Value = 132.2322999999997
$ws.Range("N138").Value = -26109.5936

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 19716.684
$ws.Range("I32").Value = 19473.88
$ws.Range("J32").Value = 21524.223
$ws.Range("K32").Value = 19473.88
$ws.Range("L32").Value = 21524.223
$ws.Range("M32").Value = -19186.88
$ws.Range("N32").Value = -22098.223

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 4714.2856
$ws.Range("I88").Value = 7500
$ws.Range("J88").Value = 2625
$ws.Range("K88").Value = 7500
$ws.Range("L88").Value = 2625
$ws.Range("M88").Value = -7094
$ws.Range("N88").Value = -3437

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 4714.2856
$ws.Range("I91").Value = 7500
$ws.Range("J91").Value = 2625
$ws.Range("K91").Value = 7500
$ws.Range("L91").Value = 2625
$ws.Range("M91").Value = -6096
$ws.Range("N91").Value = -5433

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 4086.182
$ws.Range("I122").Value = 4364.3
$ws.Range("J122").Value = 3658.3076
$ws.Range("K122").Value = 13092.9
$ws.Range("L122").Value = 10974.9228
$ws.Range("M122").Value = -10642.9
$ws.Range("N122").Value = -15874.9228

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 16034.863
$ws.Range("I132").Value = 20407.018
$ws.Range("J132").Value = 3608.7368
$ws.Range("K132").Value = 61221.054
$ws.Range("L132").Value = 10826.2104
$ws.Range("M132").Value = -58691.054
$ws.Range("N132").Value = -15886.2104

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker | Titanium Ingot
$ws.Range("H80").Value = 150.6
$ws.Range("I80").Value = 118
$ws.Range("J80").Value = 172.33333
$ws.Range("K80").Value = 118
$ws.Range("L80").Value = 172.33333
$ws.Range("M80").Value = 880
$ws.Range("N80").Value = -2168.33333

# Row 83: Attack on Titanium (L) | Titanium Ingot
$ws.Range("H83").Value = 150.6
$ws.Range("I83").Value = 118
$ws.Range("J83").Value = 172.33333
$ws.Range("K83").Value = 590
$ws.Range("L83").Value = 861.6666499999999
$ws.Range("M83").Value = 4402
$ws.Range("N83").Value = -10845.66665

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2596.1
$ws.Range("I134").Value = 2158.0417
$ws.Range("J134").Value = 4348.3335
$ws.Range("K134").Value = 6474.125100000001
$ws.Range("L134").Value = 13045.0005
$ws.Range("M134").Value = -3939.125100000001
$ws.Range("N134").Value = -18115.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2549.0435
$ws.Range("I31").Value = 1925.6482
$ws.Range("J31").Value = 4793.2666
$ws.Range("K31").Value = 1925.6482
$ws.Range("L31").Value = 4793.2666
$ws.Range("M31").Value = -1630.6482
$ws.Range("N31").Value = -5383.2666

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2549.0435
$ws.Range("I34").Value = 1925.6482
$ws.Range("J34").Value = 4793.2666
$ws.Range("K34").Value = 1925.6482
$ws.Range("L34").Value = 4793.2666
$ws.Range("M34").Value = -1723.6482
$ws.Range("N34").Value = -5197.2666

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 1553
$ws.Range("I122").Value = 964.7
$ws.Range("J122").Value = 2141.3
$ws.Range("K122").Value = 2894.1
$ws.Range("L122").Value = 6423.900000000001
$ws.Range("M122").Value = -444.1000000000004
$ws.Range("N122").Value = -11323.9

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2287.8
$ws.Range("I132").Value = 1083.826
$ws.Range("J132").Value = 4595.4165
$ws.Range("K132").Value = 3251.478
$ws.Range("L132").Value = 13786.2495
$ws.Range("M132").Value = -721.4780000000001
$ws.Range("N132").Value = -18846.2495

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 1388.1273
$ws.Range("I134").Value = 959.2895
$ws.Range("J134").Value = 2346.7058
$ws.Range("K134").Value = 2877.8685
$ws.Range("L134").Value = 7040.117400000001
$ws.Range("M134").Value = -342.8685
$ws.Range("N134").Value = -12110.1174

$ws = $wb.Worksheets.Item("CUL")
# Row 98: Sweet Kiss of Death | Rice Vinegar
$ws.Range("H98").Value = 5156.9287
$ws.Range("J98").Value = 6463.364
$ws.Range("L98").Value = 19390.092
$ws.Range("N98").Value = -22386.092

# Row 129: Comfort Food | Yakow Moussaka
$ws.Range("H129").Value = 1585.8
$ws.Range("I129").Value = 910
$ws.Range("J129").Value = 1677.9546
$ws.Range("K129").Value = 2730
$ws.Range("L129").Value = 5033.8638
$ws.Range("M129").Value = 2270
$ws.Range("N129").Value = -15033.8638

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 501486.16
$ws.Range("I122").Value = 1001069.7
$ws.Range("J122").Value = 1902.6666
$ws.Range("K122").Value = 3003209.1
$ws.Range("L122").Value = 5707.9998
$ws.Range("M122").Value = -3000759.1
$ws.Range("N122").Value = -10607.9998

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3380.2856
$ws.Range("I132").Value = 3279.3784
$ws.Range("J132").Value = 3576.7896
$ws.Range("K132").Value = 9838.135200000001
$ws.Range("L132").Value = 10730.3688
$ws.Range("M132").Value = -7308.135200000001
$ws.Range("N132").Value = -15790.3688

# Row 136: Shiny and Good | Pink Beryl
$ws.Range("H136").Value = 27500
$ws.Range("J136").Value = 27500
$ws.Range("L136").Value = 82500
$ws.Range("N136").Value = -87600

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 3023.8262
$ws.Range("I7").Value = 3061.4614
$ws.Range("J7").Value = 2974.9
$ws.Range("K7").Value = 3061.4614
$ws.Range("L7").Value = 2974.9
$ws.Range("M7").Value = -2949.4614
$ws.Range("N7").Value = -3198.9

# Row 112: A Slippery Slope | Gliderskin Boots of Casting
$ws.Range("H112").Value = 33432.547
$ws.Range("J112").Value = 33432.547
$ws.Range("L112").Value = 33432.547
$ws.Range("N112").Value = -36386.547

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 3023.8262
$ws.Range("I126").Value = 3061.4614
$ws.Range("J126").Value = 2974.9
$ws.Range("K126").Value = 9184.3842
$ws.Range("L126").Value = 8924.700000000001
$ws.Range("M126").Value = -6714.3842
$ws.Range("N126").Value = -13864.7

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 3528.3147
$ws.Range("I136").Value = 1927.4634
$ws.Range("J136").Value = 8577.154
$ws.Range("K136").Value = 5782.3902
$ws.Range("L136").Value = 25731.462
$ws.Range("M136").Value = -3232.3902
$ws.Range("N136").Value = -30831.462

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 40041400
$ws.Range("I122").Value = 52685530
$ws.Range("J122").Value = 1652.5
$ws.Range("K122").Value = 158056590
$ws.Range("L122").Value = 4957.5
$ws.Range("M122").Value = -158054140
$ws.Range("N122").Value = -9857.5

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1491
$ws.Range("I126").Value = 878.9091
$ws.Range("J126").Value = 2332.625
$ws.Range("K126").Value = 2636.7273
$ws.Range("L126").Value = 6997.875
$ws.Range("M126").Value = -166.7273
$ws.Range("N126").Value = -11937.875

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2086.8438
$ws.Range("I132").Value = 1280.0952
$ws.Range("J132").Value = 3627
$ws.Range("K132").Value = 3840.2856
$ws.Range("L132").Value = 10881
$ws.Range("M132").Value = -1310.2856
$ws.Range("N132").Value = -15941

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 15319924
$ws.Range("I136").Value = 20855504
$ws.Range("J136").Value = 558377.75
$ws.Range("K136").Value = 62566512
$ws.Range("L136").Value = 1675133.25
$ws.Range("M136").Value = -62563962
$ws.Range("N136").Value = -1680233.25

# Row 137: Traditional Trousers | Sarcenet Slops of Aiming
$ws.Range("H137").Value = 39789.168
$ws.Range("I137").Value = 32650
$ws.Range("J137").Value = 41217
$ws.Range("K137").Value = 32650
$ws.Range("L137").Value = 41217
$ws.Range("M137").Value = -27550
$ws.Range("N137").Value = -51417
